# Update the manualStatus column (I) from a raw numeric 128 to a
# manually-entered status string "[128]", for every data row (2-25).
# This also mirrors the row-height side effect (15 -> 13.8) that Excel
# applied to those rows, and moves the active selection to I25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 25

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "[128]"
    $ws.Rows.Item($r).RowHeight = 13.8
}

$ws.Range("I25").Select()
